$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 84, pushing the existing rows 84-93 down to 85-94.
$ws.Rows(84).Insert()

# Populate the newly inserted row 84 with the new record's data.
$ws.Cells.Item(84, 1).Value = 10
$ws.Cells.Item(84, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(84, 3).Value = "La Araucanía"
$ws.Cells.Item(84, 4).Value = 45212
$ws.Cells.Item(84, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(84, 5).Value = 9
$ws.Cells.Item(84, 6).Value = 100112026
$ws.Cells.Item(84, 7).Value = "Haba"
$ws.Cells.Item(84, 8).Value = "Sin especificar"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 10
$ws.Cells.Item(84, 11).Value = 14000
$ws.Cells.Item(84, 12).Value = 14000
$ws.Cells.Item(84, 13).Value = 14000
$ws.Cells.Item(84, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(84, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(84, 16).Value = 560
$ws.Cells.Item(84, 17).Value = 25
$ws.Cells.Item(84, 18).Value = "Hortaliza"
